$p = $ppt.ActivePresentation
try {
    $d2 = $p.Designs.Add("testTheme")
    Write-Output ("Added: " + $d2.Name)
} catch {
    Write-Output ("error on add: " + $_.Exception.Message)
}
Write-Output ("Designs.Count after: " + $p.Designs.Count)
